$d = $word.ActiveDocument

$d.Content.Find.Execute("87×41=3567", $true, $false, $false, $false, $false, $true, 1, $false, "62×88=5456", 2) | Out-Null
$d.Content.Find.Execute("94×95=8930", $true, $false, $false, $false, $false, $true, 1, $false, "56×56=3136", 2) | Out-Null
$d.Content.Find.Execute("97×53=5141", $true, $false, $false, $false, $false, $true, 1, $false, "92×87=8004", 2) | Out-Null
$d.Content.Find.Execute("40×16=640", $true, $false, $false, $false, $false, $true, 1, $false, "41×56=2296", 2) | Out-Null
$d.Content.Find.Execute("91×46=4186", $true, $false, $false, $false, $false, $true, 1, $false, "38×33=1254", 2) | Out-Null
$d.Content.Find.Execute("78×58=4524", $true, $false, $false, $false, $false, $true, 1, $false, "23×20=460", 2) | Out-Null
$d.Content.Find.Execute("21×75=1575", $true, $false, $false, $false, $false, $true, 1, $false, "93×53=4929", 2) | Out-Null
$d.Content.Find.Execute("32×44=1408", $true, $false, $false, $false, $false, $true, 1, $false, "47×89=4183", 2) | Out-Null
$d.Content.Find.Execute("69×25=1725", $true, $false, $false, $false, $false, $true, 1, $false, "99×68=6732", 2) | Out-Null
$d.Content.Find.Execute("54×76=4104", $true, $false, $false, $false, $false, $true, 1, $false, "97×29=2813", 2) | Out-Null
$d.Content.Find.Execute("96×14=1344", $true, $false, $false, $false, $false, $true, 1, $false, "97×33=3201", 2) | Out-Null
$d.Content.Find.Execute("73×37=2701", $true, $false, $false, $false, $false, $true, 1, $false, "50×90=4500", 2) | Out-Null
$d.Content.Find.Execute("33×99=3267", $true, $false, $false, $false, $false, $true, 1, $false, "28×44=1232", 2) | Out-Null
$d.Content.Find.Execute("42×78=3276", $true, $false, $false, $false, $false, $true, 1, $false, "45×25=1125", 2) | Out-Null
$d.Content.Find.Execute("76×71=5396", $true, $false, $false, $false, $false, $true, 1, $false, "85×19=1615", 2) | Out-Null
$d.Content.Find.Execute("25×88=2200", $true, $false, $false, $false, $false, $true, 1, $false, "21×45=945", 2) | Out-Null
$d.Content.Find.Execute("96×33=3168", $true, $false, $false, $false, $false, $true, 1, $false, "27×71=1917", 2) | Out-Null
$d.Content.Find.Execute("73×24=1752", $true, $false, $false, $false, $false, $true, 1, $false, "60×51=3060", 2) | Out-Null
$d.Content.Find.Execute("52×73=3796", $true, $false, $false, $false, $false, $true, 1, $false, "44×50=2200", 2) | Out-Null
$d.Content.Find.Execute("85×65=5525", $true, $false, $false, $false, $false, $true, 1, $false, "33×21=693", 2) | Out-Null
$d.Content.Find.Execute("88×68=5984", $true, $false, $false, $false, $false, $true, 1, $false, "18×88=1584", 2) | Out-Null
$d.Content.Find.Execute("68×21=1428", $true, $false, $false, $false, $false, $true, 1, $false, "74×15=1110", 2) | Out-Null
$d.Content.Find.Execute("36×15=540", $true, $false, $false, $false, $false, $true, 1, $false, "14×77=1078", 2) | Out-Null
$d.Content.Find.Execute("41×97=3977", $true, $false, $false, $false, $false, $true, 1, $false, "16×89=1424", 2) | Out-Null
$d.Content.Find.Execute("34×31=1054", $true, $false, $false, $false, $false, $true, 1, $false, "14×82=1148", 2) | Out-Null
